$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 714308.4399999999
$ws.Range("I6").Value = 833343.2
$ws.Range("K6").Value = 2500029.6
$ws.Range("M6").Value = -2499917.6
$ws.Range("H9").Value = 97.5
$ws.Range("I9").Value = 97
$ws.Range("J9").Value = 98
$ws.Range("K9").Value = 97
$ws.Range("L9").Value = 98
$ws.Range("M9").Value = 72
$ws.Range("N9").Value = -436
$ws.Range("H11").Value = 35.333332
$ws.Range("I11").Value = 35.333332
$ws.Range("K11").Value = 35.333332
$ws.Range("M11").Value = 104.666668
$ws.Range("H12").Value = 199.5
$ws.Range("I12").Value = 199.5
$ws.Range("K12").Value = 199.5
$ws.Range("M12").Value = -29.5
$ws.Range("H21").Value = 3849.5
$ws.Range("J21").Value = 3849.5
$ws.Range("L21").Value = 3849.5
$ws.Range("N21").Value = -4785.5
$ws.Range("H23").Value = 3849.5
$ws.Range("J23").Value = 3849.5
$ws.Range("L23").Value = 3849.5
$ws.Range("N23").Value = -4317.5
$ws.Range("H39").Value = 474.66666
$ws.Range("I39").Value = 474.66666
$ws.Range("K39").Value = 1423.99998
$ws.Range("M39").Value = -1127.99998
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H51").Value = 5002
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 5002
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 5002
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -5970
$ws.Range("H62").Value = 3499.5
$ws.Range("I62").Value = 999
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 999
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -375
$ws.Range("N62").Value = -7248
$ws.Range("H65").Value = 3499.5
$ws.Range("I65").Value = 999
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 4995
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -1875
$ws.Range("N65").Value = -36240
$ws.Range("H69").Value = 7378.8335
$ws.Range("J69").Value = 7378.8335
$ws.Range("L69").Value = 22136.5005
$ws.Range("N69").Value = -23884.5005
$ws.Range("H72").Value = 7378.8335
$ws.Range("J72").Value = 7378.8335
$ws.Range("L72").Value = 66409.5015
$ws.Range("N72").Value = -75145.5015
$ws.Range("H76").Value = 5997
$ws.Range("I76").Value = 5997
$ws.Range("K76").Value = 5997
$ws.Range("M76").Value = -5682
$ws.Range("H79").Value = 5997
$ws.Range("I79").Value = 5997
$ws.Range("K79").Value = 5997
$ws.Range("M79").Value = -4905
$ws.Range("H80").Value = 1926.6666
$ws.Range("I80").Value = 1780
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 5340
$ws.Range("L80").Value = 6000
$ws.Range("M80").Value = -4342
$ws.Range("N80").Value = -7996
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H83").Value = 1926.6666
$ws.Range("I83").Value = 1780
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 16020
$ws.Range("L83").Value = 18000
$ws.Range("M83").Value = -11028
$ws.Range("N83").Value = -27984
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H92").Value = 3051.4443
$ws.Range("I92").Value = 2183
$ws.Range("K92").Value = 2183
$ws.Range("M92").Value = -935
$ws.Range("H96").Value = 812.1429000000001
$ws.Range("I96").Value = 655
$ws.Range("J96").Value = 1021.6667
$ws.Range("K96").Value = 1965
$ws.Range("L96").Value = 3065.0001
$ws.Range("M96").Value = -592
$ws.Range("N96").Value = -5811.0001
$ws.Range("H98").Value = 13312.125
$ws.Range("J98").Value = 14666.667
$ws.Range("L98").Value = 14666.667
$ws.Range("N98").Value = -17662.667
$ws.Range("H122").Value = 13312.125
$ws.Range("J122").Value = 14666.667
$ws.Range("L122").Value = 44000.001
$ws.Range("N122").Value = -48900.001
$ws.Range("H125").Value = 2329.6667
$ws.Range("I125").Value = 1994.5
$ws.Range("J125").Value = 3000
$ws.Range("K125").Value = 17950.5
$ws.Range("L125").Value = 27000
$ws.Range("M125").Value = -15490.5
$ws.Range("N125").Value = -31920
$ws.Range("H138").Value = 3253.0476
$ws.Range("J138").Value = 3384.4
$ws.Range("L138").Value = 10153.2
$ws.Range("N138").Value = -20433.2

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8000
$ws.Range("I32").Value = 8000
$ws.Range("K32").Value = 8000
$ws.Range("M32").Value = -7713
$ws.Range("H74").Value = 5221.778
$ws.Range("J74").Value = 7800
$ws.Range("L74").Value = 7800
$ws.Range("N74").Value = -9548
$ws.Range("H77").Value = 5221.778
$ws.Range("J77").Value = 7800
$ws.Range("L77").Value = 39000
$ws.Range("N77").Value = -47736

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 303.66666
$ws.Range("I86").Value = 303.66666
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 303.66666
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 819.33334
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 303.66666
$ws.Range("I89").Value = 303.66666
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 1518.3333
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 4097.6667
$ws.Range("N89").ClearContents()
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H134").Value = 7519.5557
$ws.Range("I134").Value = 4002
$ws.Range("J134").Value = 11916.5
$ws.Range("K134").Value = 12006
$ws.Range("L134").Value = 35749.5
$ws.Range("M134").Value = -9471
$ws.Range("N134").Value = -40819.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 1849.5555
$ws.Range("J32").Value = 4166.6665
$ws.Range("L32").Value = 4166.6665
$ws.Range("N32").Value = -4798.6665
$ws.Range("H141").Value = 874721.5
$ws.Range("J141").Value = 1146295.6
$ws.Range("L141").Value = 1146295.6
$ws.Range("N141").Value = -1156655.6

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 949.8946999999999
$ws.Range("J2").Value = 3900
$ws.Range("L2").Value = 23400
$ws.Range("N2").Value = -23626
$ws.Range("H10").Value = 155.41667
$ws.Range("I10").Value = 49.166668
$ws.Range("J10").Value = 261.66666
$ws.Range("K10").Value = 147.500004
$ws.Range("L10").Value = 784.9999799999999
$ws.Range("M10").Value = -8.50000399999999
$ws.Range("N10").Value = -1062.99998
$ws.Range("H57").Value = 97.5
$ws.Range("I57").Value = 97.5
$ws.Range("K57").Value = 292.5
$ws.Range("M57").Value = 266.5
$ws.Range("H58").Value = 1000
$ws.Range("J58").Value = 1000
$ws.Range("L58").Value = 3000
$ws.Range("N58").Value = -3256
$ws.Range("H81").Value = 2995
$ws.Range("J81").Value = 2995
$ws.Range("L81").Value = 8985
$ws.Range("N81").Value = -11231
$ws.Range("H84").Value = 2995
$ws.Range("J84").Value = 2995
$ws.Range("L84").Value = 26955
$ws.Range("N84").Value = -38187
$ws.Range("H140").Value = 1162.25
$ws.Range("I140").Value = 1162.25
$ws.Range("K140").Value = 3486.75
$ws.Range("M140").Value = 1693.25

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5000
$ws.Range("I70").Value = 5000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 5000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -4730
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 5000
$ws.Range("I73").Value = 5000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 5000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -4064
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 1273.25
$ws.Range("J80").Value = 2350
$ws.Range("L80").Value = 2350
$ws.Range("N80").Value = -4346
$ws.Range("H83").Value = 1273.25
$ws.Range("J83").Value = 2350
$ws.Range("L83").Value = 11750
$ws.Range("N83").Value = -21734
$ws.Range("H113").Value = 500000000
$ws.Range("I113").Value = 500000000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 500000000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -499997830
$ws.Range("N113").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 11036.75
$ws.Range("I58").Value = 9101.666999999999
$ws.Range("J58").Value = 16842
$ws.Range("K58").Value = 9101.666999999999
$ws.Range("L58").Value = 16842
$ws.Range("M58").Value = -8841.666999999999
$ws.Range("N58").Value = -17362
$ws.Range("H74").Value = 36039.4
$ws.Range("J74").Value = 25000
$ws.Range("L74").Value = 25000
$ws.Range("N74").Value = -26996
$ws.Range("H77").Value = 36039.4
$ws.Range("J77").Value = 25000
$ws.Range("L77").Value = 75000
$ws.Range("N77").Value = -84984
$ws.Range("H132").Value = 12857.143
$ws.Range("I132").Value = 10000
$ws.Range("J132").Value = 20000
$ws.Range("K132").Value = 30000
$ws.Range("L132").Value = 60000
$ws.Range("M132").Value = -27470
$ws.Range("N132").Value = -65060

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 16000.333
$ws.Range("I136").Value = 14001
$ws.Range("J136").Value = 17000
$ws.Range("K136").Value = 42003
$ws.Range("L136").Value = 51000
$ws.Range("M136").Value = -39453
$ws.Range("N136").Value = -56100
